$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Gene symbol corrections
$ws.Range("C11").Value = "IL12A"
$ws.Range("C30").Value = "LTA"

# Move the active selection to C31 (just below the data, matching the
# author's WIP cursor position while comparing liver-quality gene sets)
$ws.Range("C31").Select()
